$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("BJ1").Value = "31-ago"

$ws.Range("BJ2").Value = 0
$ws.Range("BJ3").Value = 12.498609596300213
$ws.Range("BJ4").Value = 14.266603824463326
$ws.Range("BJ5").Value = 12.896729109586072
$ws.Range("BJ6").Value = 0
$ws.Range("BJ7").Value = 5.231148683730761
$ws.Range("BJ8").Value = 10.093403248029892
$ws.Range("BJ9").Value = 4.3657291897079604
$ws.Range("BJ10").Value = 24.37233824545234
$ws.Range("BJ11").Value = 15.690741044455448
$ws.Range("BJ12").Value = 0
$ws.Range("BJ13").Value = 8.2846987450042402
$ws.Range("BJ14").Value = 0
$ws.Range("BJ15").Value = 0
$ws.Range("BJ16").Value = 17.861450191821788
$ws.Range("BJ17").Value = 0
$ws.Range("BJ18").Value = 0

$ws.Range("BK4").Select()
